$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, D (date serial), I (Calidad), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), N (Unidad comercializacion),
# P (Precio $/Kg), Q (Kg o Unidades)
$rowData = @(
    @(60, 44942, 'Primera', 520, 14000, 15000, 14500, '$/caja 13 kilos', 1115, 13),
    @(61, 44669, 'Primera', 610, 10000, 11000, 10500, '$/caja 13 kilos', 808, 13),
    @(62, 44424, 'Primera', 700, 13000, 14000, 13500, '$/caja 13 kilos', 1038, 13),
    @(63, 44424, 'Segunda', 430, 12000, 12000, 12000, '$/caja 13 kilos', 923, 13),
    @(64, 44459, 'Primera', 970, 13000, 14000, 13495, '$/caja 13 kilos', 1038, 13),
    @(65, 44459, 'Segunda', 520, 11000, 12000, 11500, '$/caja 13 kilos', 885, 13),
    @(66, 44263, 'Primera', 250, 40000, 40000, 40000, '$/caja 15 kilos', 2667, 15),
    @(67, 44550, 'Primera', 790, 11000, 12000, 11494, '$/caja 13 kilos', 884, 13),
    @(68, 44550, 'Segunda', 430, 9000, 10000, 9500, '$/caja 13 kilos', 731, 13),
    @(69, 44319, 'Primera', 340, 24000, 25000, 24500, '$/caja 13 kilos', 1885, 13),
    @(70, 44319, 'Segunda', 160, 22000, 22000, 22000, '$/caja 13 kilos', 1692, 13),
    @(71, 44921, 'Primera', 610, 15000, 15000, 15000, '$/caja 13 kilos', 1154, 13),
    @(72, 44900, 'Primera', 700, 15000, 17000, 16000, '$/caja 13 kilos', 1231, 13),
    @(73, 44445, 'Primera', 790, 13000, 14000, 13494, '$/caja 13 kilos', 1038, 13),
    @(74, 44445, 'Segunda', 340, 11000, 12000, 11500, '$/caja 13 kilos', 885, 13),
    @(75, 44914, 'Primera', 520, 13000, 14000, 13500, '$/caja 13 kilos', 1038, 13),
    @(76, 44606, 'Primera', 520, 17000, 18000, 17500, '$/caja 13 kilos', 1346, 13),
    @(77, 44571, 'Primera', 610, 12000, 13000, 12500, '$/caja 13 kilos', 962, 13),
    @(78, 44571, 'Segunda', 106, 10000, 10000, 10000, '$/caja 13 kilos', 769, 13),
    @(79, 44396, 'Primera', 770, 17000, 18000, 17494, '$/caja 13 kilos', 1346, 13),
    @(80, 44396, 'Segunda', 340, 16000, 16000, 16000, '$/caja 13 kilos', 1231, 13),
    @(81, 44809, 'Primera', 350, 13000, 15000, 14143, '$/caja 13 kilos', 1088, 13),
    @(82, 44809, 'Segunda', 160, 10000, 10000, 10000, '$/caja 13 kilos', 769, 13),
    @(83, 44816, 'Primera', 790, 14000, 14000, 14000, '$/caja 13 kilos', 1077, 13),
    @(84, 44613, 'Primera', 790, 16000, 17000, 16494, '$/caja 13 kilos', 1269, 13),
    @(85, 44371, 'Primera', 160, 20000, 21000, 20500, '$/caja 13 kilos', 1577, 13),
    @(86, 44340, 'Primera', 250, 20000, 20000, 20000, '$/caja 13 kilos', 1538, 13),
    @(87, 44340, 'Segunda', 160, 18000, 18000, 18000, '$/caja 13 kilos', 1385, 13),
    @(88, 44200, 'Primera', 520, 30000, 30000, 30000, '$/caja 13 kilos', 2308, 13),
    @(89, 44200, 'Segunda', 340, 25000, 25000, 25000, '$/caja 13 kilos', 1923, 13),
    @(90, 44641, 'Primera', 610, 14000, 15000, 14500, '$/caja 13 kilos', 1115, 13),
    @(91, 44221, 'Segunda', 180, 35000, 35000, 35000, '$/caja 13 kilos', 2692, 13),
    @(92, 44305, 'Primera', 340, 24000, 24000, 24000, '$/caja 13 kilos', 1846, 13),
    @(93, 44305, 'Segunda', 160, 20000, 20000, 20000, '$/caja 13 kilos', 1538, 13),
    @(94, 44494, 'Primera', 780, 15000, 15000, 15000, '$/caja 13 kilos', 1154, 13),
    @(95, 44760, 'Primera', 430, 16000, 17000, 16500, '$/caja 13 kilos', 1269, 13),
    @(96, 44592, 'Primera', 400, 9000, 10000, 9575, '$/caja 13 kilos', 737, 13),
    @(97, 44592, 'Segunda', 100, 8000, 8000, 8000, '$/caja 13 kilos', 615, 13),
    @(98, 44858, 'Primera', 440, 13000, 14000, 13545, '$/caja 13 kilos', 1042, 13),
    @(99, 44725, 'Primera', 610, 14000, 15000, 14500, '$/caja 13 kilos', 1115, 13),
    @(100, 44298, 'Primera', 340, 24000, 25000, 24500, '$/caja 13 kilos', 1885, 13),
    @(101, 44627, 'Primera', 790, 14000, 15000, 14494, '$/caja 13 kilos', 1115, 13),
    @(102, 44627, 'Segunda', 340, 13000, 13000, 13000, '$/caja 13 kilos', 1000, 13),
    @(103, 44487, 'Primera', 1150, 14000, 15000, 14500, '$/caja 13 kilos', 1115, 13),
    @(104, 44487, 'Segunda', 610, 12000, 12000, 12000, '$/caja 13 kilos', 923, 13),
    @(105, 44830, 'Primera', 450, 14000, 15000, 14556, '$/caja 13 kilos', 1120, 13),
    @(106, 44403, 'Primera', 700, 16000, 17000, 16500, '$/caja 13 kilos', 1269, 13),
    @(107, 44403, 'Segunda', 430, 15000, 15000, 15000, '$/caja 13 kilos', 1154, 13),
    @(108, 44277, 'Primera', 250, 38000, 38000, 38000, '$/caja 13 kilos', 2923, 13),
    @(109, 44277, 'Segunda', 160, 35000, 35000, 35000, '$/caja 13 kilos', 2692, 13),
    @(110, 44585, 'Primera', 790, 10000, 11000, 10494, '$/caja 13 kilos', 807, 13),
    @(111, 44585, 'Segunda', 340, 9000, 9000, 9000, '$/caja 13 kilos', 692, 13),
    @(112, 44893, 'Primera', 610, 15000, 17000, 16000, '$/caja 13 kilos', 1231, 13),
    @(113, 44382, 'Primera', 790, 14000, 15000, 14500, '$/caja 13 kilos', 1115, 13),
    @(114, 44382, 'Segunda', 430, 12000, 12000, 12000, '$/caja 13 kilos', 923, 13),
    @(115, 44879, 'Primera', 380, 14000, 15000, 14526, '$/caja 13 kilos', 1117, 13),
    @(116, 44781, 'Primera', 430, 15000, 16000, 15500, '$/caja 13 kilos', 1192, 13),
    @(117, 44249, 'Primera', 250, 39000, 42000, 40500, '$/caja 13 kilos', 3115, 13),
    @(118, 44690, 'Primera', 790, 12000, 13000, 12494, '$/caja 13 kilos', 961, 13)
)

foreach ($r in $rowData) {
    $row = $r[0]
    $ws.Cells.Item($row, 4).Value = $r[1]          # D - Fecha
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 9).Value = $r[2]          # I - Calidad
    $ws.Cells.Item($row, 10).Value = $r[3]         # J - Volumen
    $ws.Cells.Item($row, 11).Value = $r[4]         # K - Precio minimo
    $ws.Cells.Item($row, 12).Value = $r[5]         # L - Precio maximo
    $ws.Cells.Item($row, 13).Value = $r[6]         # M - Precio promedio ponderado
    $ws.Cells.Item($row, 14).Value = $r[7]         # N - Unidad de comercializacion
    $ws.Cells.Item($row, 16).Value = $r[8]         # P - Precio $/Kg
    $ws.Cells.Item($row, 17).Value = $r[9]         # Q - Kg o Unidades
}

# Row 118 is brand new; fill in the constant columns (A, B, C, E, F, G, H, O, R)
# that are identical across every data row in this sheet.
$ws.Cells.Item(118, 1).Value = 9
$ws.Cells.Item(118, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(118, 3).Value = "Metropolitana"
$ws.Cells.Item(118, 5).Value = 13
$ws.Cells.Item(118, 6).Value = 100114007
$ws.Cells.Item(118, 7).Value = "Jengibre"
$ws.Cells.Item(118, 8).Value = "Sin especificar"
$ws.Cells.Item(118, 15).Value = "Perú"
$ws.Cells.Item(118, 18).Value = "Hortaliza"
